# Fix exam security loopholes and student re-entry issues
#
# The template's "username" column (H) stored a bare login handle
# (johndoe001). Students were re-using that same handle to re-enter the
# exam, so it is replaced with the student's actual NMIMS e-mail address
# and turned into a clickable mailto: hyperlink (matching the Hyperlink
# cell style Excel applies automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("H2")

# Swap the plain username for the student's e-mail address.
$cell.Value = "john.doe01@nmims.in"

# Turn it into a real mailto: hyperlink (adds the Hyperlink font/style and
# the worksheet-level hyperlink relationship, like Excel does natively).
$ws.Hyperlinks.Add($cell, "mailto:john.doe01@nmims.in", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "john.doe01@nmims.in")
